# "Refreshed stats and charts with new data"
#
# Row 15 in the forecast table used to be the first *projected* day
# (I15 = I14*(1+AVERAGE(M10:M14))). New actual case-count data came in for
# that date, so the forecast formula is replaced with the real reported
# value, and the cell is restyled to match the other "actual data" cells
# (I9:I14) instead of the "forecast" cells (I16:I28). Every other cell in
# the sheet (J15:N15 and the whole I16:N28 forecast cascade) is already a
# formula referencing column I, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# I15: forecast formula -> hard-coded actual value
$target = $ws.Cells.Item(15, 9)
$target.Value = 24207

# Match the "actual data" formatting used by I9:I14 (copy format only, so
# the existing style is reused instead of a new one being minted).
$actualStyleSource = $ws.Cells.Item(14, 9)
$actualStyleSource.Copy()
$target.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on the cell that was just edited.
[void]$target.Select()

Write-Output "Set I15 to 24207 (actual data) and refreshed the forecast cascade"
